$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 60000
$ws.Range("J81").Value = 60000
$ws.Range("L81").Value = 60000
$ws.Range("N81").Value = -61996
$ws.Range("H84").Value = 60000
$ws.Range("J84").Value = 60000
$ws.Range("L84").Value = 180000
$ws.Range("N84").Value = -189984
$ws.Range("H87").Value = 36997.332
$ws.Range("J87").Value = 36997.332
$ws.Range("L87").Value = 36997.332
$ws.Range("N87").Value = -39493.332
$ws.Range("H90").Value = 36997.332
$ws.Range("J90").Value = 36997.332
$ws.Range("L90").Value = 110991.996
$ws.Range("N90").Value = -123471.996
$ws.Range("H120").Value = 49716.668
$ws.Range("J120").Value = 49716.668
$ws.Range("L120").Value = 49716.668
$ws.Range("N120").Value = -59392.668
$ws.Range("H137").Value = 4330.244
$ws.Range("I137").Value = 1382.75
$ws.Range("K137").Value = 4148.25
$ws.Range("M137").Value = -1598.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 49116.57
$ws.Range("J80").Value = 49116.57
$ws.Range("L80").Value = 49116.57
$ws.Range("N80").Value = -51112.57
$ws.Range("H83").Value = 49116.57
$ws.Range("J83").Value = 49116.57
$ws.Range("L83").Value = 147349.71
$ws.Range("N83").Value = -157333.71
$ws.Range("H121").Value = 45238.332
$ws.Range("J121").Value = 45238.332
$ws.Range("L121").Value = 45238.332
$ws.Range("N121").Value = -48732.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3239.8062
$ws.Range("I31").Value = 1148
$ws.Range("J31").Value = 3918.2297
$ws.Range("K31").Value = 1148
$ws.Range("L31").Value = 3918.2297
$ws.Range("M31").Value = -853
$ws.Range("N31").Value = -4508.2297
$ws.Range("H34").Value = 3239.8062
$ws.Range("I34").Value = 1148
$ws.Range("J34").Value = 3918.2297
$ws.Range("K34").Value = 1148
$ws.Range("L34").Value = 3918.2297
$ws.Range("M34").Value = -946
$ws.Range("N34").Value = -4322.2297
$ws.Range("H100").Value = 47675
$ws.Range("J100").Value = 47675
$ws.Range("L100").Value = 47675
$ws.Range("N100").Value = -49839
$ws.Range("H124").Value = 45318
$ws.Range("J124").Value = 45318
$ws.Range("L124").Value = 45318
$ws.Range("N124").Value = -50228
$ws.Range("H125").Value = 38659
$ws.Range("J125").Value = 38659
$ws.Range("L125").Value = 38659
$ws.Range("N125").Value = -43579

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 38310
$ws.Range("J118").Value = 38310
$ws.Range("L118").Value = 38310
$ws.Range("N118").Value = -41624
$ws.Range("H120").Value = 37992
$ws.Range("J120").Value = 37992
$ws.Range("L120").Value = 37992
$ws.Range("N120").Value = -47668
$ws.Range("H125").Value = 34814.668
$ws.Range("J125").Value = 34814.668
$ws.Range("L125").Value = 34814.668
$ws.Range("N125").Value = -39734.668
$ws.Range("H127").Value = 47303
$ws.Range("J127").Value = 47303
$ws.Range("L127").Value = 47303
$ws.Range("N127").Value = -57223
$ws.Range("H131").Value = 42326
$ws.Range("J131").Value = 42326
$ws.Range("L131").Value = 42326
$ws.Range("N131").Value = -52406

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 38374
$ws.Range("J92").Value = 38374
$ws.Range("L92").Value = 38374
$ws.Range("N92").Value = -43366
$ws.Range("H99").Value = 27471
$ws.Range("I99").Value = 18309.75
$ws.Range("J99").Value = 34800
$ws.Range("K99").Value = 18309.75
$ws.Range("L99").Value = 34800
$ws.Range("M99").Value = -15314.75
$ws.Range("N99").Value = -40790
$ws.Range("H109").Value = 35281
$ws.Range("J109").Value = 35281
$ws.Range("L109").Value = 35281
$ws.Range("N109").Value = -38055
$ws.Range("H117").Value = 39058.668
$ws.Range("J117").Value = 39058.668
$ws.Range("L117").Value = 39058.668
$ws.Range("N117").Value = -48236.668
$ws.Range("H123").Value = 26660.5
$ws.Range("J123").Value = 32880.668
$ws.Range("L123").Value = 32880.668
$ws.Range("N123").Value = -42680.668
$ws.Range("H125").Value = 49711
$ws.Range("J125").Value = 49711
$ws.Range("L125").Value = 49711
$ws.Range("N125").Value = -59551
$ws.Range("H129").Value = 44425
$ws.Range("J129").Value = 44425
$ws.Range("L129").Value = 44425
$ws.Range("N129").Value = -54425
$ws.Range("H131").Value = 43326
$ws.Range("J131").Value = 43326
$ws.Range("L131").Value = 43326
$ws.Range("N131").Value = -53406
$ws.Range("H137").Value = 29966.666
$ws.Range("J137").Value = 29966.666
$ws.Range("L137").Value = 29966.666
$ws.Range("N137").Value = -40166.666
$ws.Range("H139").Value = 71999.8
$ws.Range("I139").Value = 160000
$ws.Range("J139").Value = 49999.75
$ws.Range("K139").Value = 160000
$ws.Range("L139").Value = 49999.75
$ws.Range("M139").Value = -154860
$ws.Range("N139").Value = -60279.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 22994
$ws.Range("I27").Value = 15500
$ws.Range("J27").Value = 30488
$ws.Range("K27").Value = 15500
$ws.Range("L27").Value = 30488
$ws.Range("M27").Value = -15431
$ws.Range("N27").Value = -30626
$ws.Range("H93").Value = 35469.855
$ws.Range("J93").Value = 35469.855
$ws.Range("L93").Value = 35469.855
$ws.Range("N93").Value = -40461.855
$ws.Range("H102").Value = 41337
$ws.Range("J102").Value = 41337
$ws.Range("L102").Value = 41337
$ws.Range("N102").Value = -47827
$ws.Range("H109").Value = 32365
$ws.Range("J109").Value = 32365
$ws.Range("L109").Value = 32365
$ws.Range("N109").Value = -35139
$ws.Range("H115").Value = 38377
$ws.Range("J115").Value = 38377
$ws.Range("L115").Value = 38377
$ws.Range("N115").Value = -41511
$ws.Range("H118").Value = 33258.668
$ws.Range("J118").Value = 42388
$ws.Range("L118").Value = 42388
$ws.Range("N118").Value = -45702
$ws.Range("H127").Value = 16284.571
$ws.Range("J127").Value = 41992
$ws.Range("L127").Value = 41992
$ws.Range("N127").Value = -51912
